$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column G, rows 4-5 (95d98036 file)
$wsOverview.Range("G4").Value = "2016-08-13 16:23:04"
$wsOverview.Range("G5").Value = "2016-08-13 16:23:04"

# zh-cn sheet: Priority (E), Correspond Handoff Datetime (H), Correspond Handback DateTime (K)
$wsZhCn.Range("E4").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"
$wsZhCn.Range("H4").Value = "2016-08-13 16:22:54"
$wsZhCn.Range("H5").Value = "2016-08-13 16:22:54"
$wsZhCn.Range("K4").Value = "2016-08-13 16:23:26"
$wsZhCn.Range("K5").Value = "2016-08-13 16:23:26"

# de-de sheet: Priority (E), Correspond Handoff Datetime (H), Correspond Handback DateTime (K)
$wsDeDe.Range("E4").Value = "mt"
$wsDeDe.Range("E5").Value = "mt"
$wsDeDe.Range("H4").Value = "2016-08-13 16:23:04"
$wsDeDe.Range("H5").Value = "2016-08-13 16:23:04"
$wsDeDe.Range("K4").Value = "2016-08-13 16:23:36"
$wsDeDe.Range("K5").Value = "2016-08-13 16:23:36"
